$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '71.345.07'
$ws.Cells.Item(2, 5).Value = '  +0.67%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.820.80'
$ws.Cells.Item(3, 5).Value = '  -0.80%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.07%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '708.77'
$ws.Cells.Item(5, 5).Value = '  +1.40%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '172.22'
$ws.Cells.Item(6, 5).Value = '  -0.12%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.819.91'
$ws.Cells.Item(7, 5).Value = '  -0.75%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.04%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.05%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +0.05%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '7.65'
$ws.Cells.Item(11, 5).Value = '  +6.25%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.91%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000255'
$ws.Cells.Item(13, 5).Value = '  -0.93%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '36.08'
$ws.Cells.Item(14, 5).Value = '  -0.32%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.464.54'
$ws.Cells.Item(15, 5).Value = '  -0.83%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.818.94'
$ws.Cells.Item(16, 5).Value = '  -1.89%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '71.279.43'
$ws.Cells.Item(17, 5).Value = '  +0.45%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +0.18%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '17.54'
$ws.Cells.Item(19, 5).Value = '  +0.75%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -0.13%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '516.96'
$ws.Cells.Item(21, 5).Value = '  +3.58%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.72'
$ws.Cells.Item(22, 5).Value = '  +0.79%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.94%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '84.74'
$ws.Cells.Item(24, 5).Value = '  -0.09%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -2.14%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '3.970.46'
$ws.Cells.Item(26, 5).Value = '  -0.87%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.07'
$ws.Cells.Item(27, 5).Value = '  -1.23%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.46'
$ws.Cells.Item(28, 5).Value = '  -1.54%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.18%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -3.30%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -2.98%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'NEARProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.43'
$ws.Cells.Item(32, 5).Value = '  -1.25%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.26'
$ws.Cells.Item(33, 5).Value = '  -0.24%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '29.20'
$ws.Cells.Item(34, 5).Value = '  -0.93%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -4.75%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '9.19'
$ws.Cells.Item(36, 5).Value = '  +0.13%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '3.790.14'
$ws.Cells.Item(37, 5).Value = '  -0.49%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.09%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -1.50%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.40%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.49%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.03'
$ws.Cells.Item(42, 5).Value = '  -2.13%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.30'
$ws.Cells.Item(43, 5).Value = '  -2.16%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.00%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.04%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '167.30'
$ws.Cells.Item(46, 5).Value = '  +2.03%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.000315'
$ws.Cells.Item(47, 5).Value = '  +1.30%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '49.35'
$ws.Cells.Item(48, 5).Value = '  +0.77%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '423.90'
$ws.Cells.Item(49, 5).Value = '  +3.01%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.38%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'ONDO'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.36'
$ws.Cells.Item(51, 5).Value = '  -1.32%  '
